# Update South Carolina overview factsheet: convert numeric "count" cells to
# text (to match formatting of the rest of the sheet), fix the Saluda County
# placeholder row, and append a "Total" row to the County sheet.

function Set-TextValue($range, [string]$text) {
    # Force the cell to Text so Excel doesn't re-parse a numeric-looking
    # string (e.g. "1,290") back into a number; then restore a plain
    # "General"/Normal style so no stray formatting is left behind.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overall sheet: A2 (No. of 990 Filers w/ Gov Grants) 1290 -> "1,290"
# ---------------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextValue $wsOverall.Range("A2") "1,290"

# ---------------------------------------------------------------------------
# County sheet
# ---------------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")

$countyCounts = @{
    2  = "5";   3  = "38";  4  = "2";   5  = "39";  6  = "6";
    7  = "2";   8  = "54";  9  = "14";  10 = "2";   11 = "192";
    12 = "17";  13 = "6";   14 = "15";  15 = "10";  16 = "8";
    17 = "14";  18 = "2";   19 = "20";  20 = "5";   21 = "5";
    22 = "41";  23 = "31";  24 = "134"; 25 = "32";  26 = "1";
    27 = "59";  28 = "7";   29 = "19";  30 = "15";  31 = "13";
    32 = "5";   33 = "39";  34 = "9";   35 = "8";   36 = "4";
    37 = "9";   38 = "11";  39 = "20";  40 = "26";  41 = "183";
    42 = "83";  43 = "23";  44 = "3";   45 = "8";   46 = "51"
}
foreach ($row in $countyCounts.Keys) {
    Set-TextValue $wsCounty.Cells.Item($row, 2) $countyCounts[$row]
}

# Row 47 (Saluda County) had placeholder zeros; replace with formatted text.
Set-TextValue $wsCounty.Cells.Item(47, 2) "0.00%"
Set-TextValue $wsCounty.Cells.Item(47, 3) "`$0"
Set-TextValue $wsCounty.Cells.Item(47, 4) "0.00%"
Set-TextValue $wsCounty.Cells.Item(47, 5) "0.00%"
Set-TextValue $wsCounty.Cells.Item(47, 6) "0.00%"

# New row 48: statewide "Total" row.
Set-TextValue $wsCounty.Cells.Item(48, 1) "Total"
Set-TextValue $wsCounty.Cells.Item(48, 2) "1,290"
Set-TextValue $wsCounty.Cells.Item(48, 3) "`$2,159,545,458"
Set-TextValue $wsCounty.Cells.Item(48, 4) "9.53%"
Set-TextValue $wsCounty.Cells.Item(48, 5) "-19.89%"
Set-TextValue $wsCounty.Cells.Item(48, 6) "70.54%"

# ---------------------------------------------------------------------------
# Congressional District sheet
# ---------------------------------------------------------------------------
$wsDistrict = $wb.Worksheets.Item("Congressional District")

$districtCounts = @{
    2 = "128"; 3 = "129"; 4 = "153"; 5 = "205"; 6 = "130"; 7 = "372"; 8 = "173"
}
foreach ($row in $districtCounts.Keys) {
    Set-TextValue $wsDistrict.Cells.Item($row, 2) $districtCounts[$row]
}
Set-TextValue $wsDistrict.Cells.Item(9, 2) "1,290"

# ---------------------------------------------------------------------------
# Size sheet
# ---------------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")

$sizeCounts = @{
    2 = "445"; 3 = "367"; 4 = "211"; 5 = "70"; 6 = "116"; 7 = "81"
}
foreach ($row in $sizeCounts.Keys) {
    Set-TextValue $wsSize.Cells.Item($row, 2) $sizeCounts[$row]
}
Set-TextValue $wsSize.Cells.Item(8, 2) "1,290"

# ---------------------------------------------------------------------------
# Subsector sheet
# ---------------------------------------------------------------------------
$wsSubsector = $wb.Worksheets.Item("Subsector")

$subsectorCounts = @{
    2  = "112"; 3  = "177"; 4  = "69";  5  = "132"; 6  = "11";
    7  = "411"; 8  = "6";   9  = "1";   10 = "115"; 11 = "32";
    12 = "210"; 13 = "14"
}
foreach ($row in $subsectorCounts.Keys) {
    Set-TextValue $wsSubsector.Cells.Item($row, 2) $subsectorCounts[$row]
}
Set-TextValue $wsSubsector.Cells.Item(14, 2) "1,290"
